$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 and B3 values
$ws.Range("B2").Value = 40425123
$ws.Range("B3").Value = 43019123

# Update D3 to new value "Soporte"
$ws.Range("D3").Value = "Soporte"

# Update column B width (widened to fit the longer ticket_id values)
$ws.Columns.Item(2).ColumnWidth = 7.83

# Update selection to A4
$ws.Range("A4").Select()
